# Updates R programs and scenarios
# Fill in row 11 (columns C:I) on both "Test 1" and "Test 2" worksheets
# with the computed scenario values, for each worksheet in the workbook.

$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 0.19508153557370633
    "D11" = -0.4640841495396
    "E11" = 0.6005884338034946
    "F11" = -0.04700000000000004
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 0.584484590860797
}

foreach ($ws in $wb.Worksheets) {
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
